$wb = $excel.ActiveWorkbook

# --- Sheet "SPN": mark some Pendente rows as Resolvido (plain value change, no highlight) ---
$wsSPN = $wb.Worksheets.Item("SPN")
$spnCells = @("I2", "I3", "I4", "I7", "I8", "I11")
foreach ($addr in $spnCells) {
    $wsSPN.Range($addr).Value = "Resolvido"
}

# --- Sheet "ITI": mark several Pendente rows as Resolvido and highlight them in yellow ---
$wsITI = $wb.Worksheets.Item("ITI")
$itiCells = @("I2", "I3", "I4", "I5", "I6", "I7", "I8", "I9", "I12", "I15", "I16", "I17", "I19", "I20", "I22", "I24", "I25", "I26", "I27", "I28", "I29")
foreach ($addr in $itiCells) {
    $cell = $wsITI.Range($addr)
    $cell.Value = "Resolvido"
    $cell.Interior.Color = 65535
}

# --- Restore the last selected cell on each sheet ---
[void]$wsSPN.Select()
[void]$wsSPN.Range("I17").Select()

[void]$wsITI.Select()
[void]$wsITI.Range("I10").Select()
